# Remove semicolon in equation area
# This script replaces semicolons used in the "variable = definition;" style
# equation/legend lines throughout the document with colons (and, for the one
# case with nothing following the semicolon, removes it outright). One
# occurrence also gains a new trailing clause (" in") instead of a colon.

$d = $word.ActiveDocument

function Replace-Unique($find, $replace) {
    $d.Content.Find.Execute($find, $false, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# 1. "...olumetric flow rate of free air;" -> "...free air" + " in" (continues
#    into "... in cubic feet per minute")
Replace-Unique "olumetric flow rate of free air;" "olumetric flow rate of free air in"

# 2. "Diameter of the leak;" -> "Diameter of the leak" (nothing follows, so the
#    semicolon is simply dropped, not replaced with a colon)
$found = $d.Content
$found.Find.Execute("Diameter of the leak", $false, $false, $false, $false, `
                     $false, $true, 1, $false, "", 0) | Out-Null
$semi = $d.Range($found.End, $found.End + 1)
if ($semi.Text -eq ";") { $semi.Text = "" }

# 3. "...mpressor inlet;" -> "...mpressor inlet:"
Replace-Unique "mpressor inlet;" "mpressor inlet:"

# 4. "...equivalent to gage line pressure;" -> "...gage line pressure:"
Replace-Unique "gage line pressure;" "gage line pressure:"

# 5. "= Inlet (atmospheric) pressure;" -> "= Inlet (atmospheric) pressure:"
Replace-Unique "= Inlet (atmospheric) pressure;" "= Inlet (atmospheric) pressure:"

# 6. " sonic volumetric flow constant;" -> " sonic volumetric flow constant:"
Replace-Unique " sonic volumetric flow constant;" " sonic volumetric flow constant:"

# 7. "= Conversion constant;" (C2, followed by "60 sec/min")
Replace-Unique "= Conversion constant; 60 sec/min" "= Conversion constant: 60 sec/min"

# 8. "...square edged orifice[footnote];" -> "...[footnote]:" (semicolon run
#    sits right after a footnote-reference run, so isolate it by position)
$found2 = $d.Content
$found2.Find.Execute("square edged orifice", $false, $false, $false, $false, `
                      $false, $true, 1, $false, "", 0) | Out-Null
$semi2 = $d.Range($found2.End + 1, $found2.End + 2)
if ($semi2.Text -eq ";") { $semi2.Text = ":" }

# 9. " constant;" (Mathematical constant / pi) -> " constant:"
Replace-Unique " constant; 3.1416" " constant: 3.1416"

# 10. " Conversion constant;" (C3, followed by "144 in") -> " Conversion constant:"
Replace-Unique " Conversion constant; 144" " Conversion constant: 144"

# 11. "= Average line temperature;" -> "= Average line temperature:"
Replace-Unique "= Average line temperature;" "= Average line temperature:"

# 12. "= Specific heat ratio of air;" -> "= Specific heat ratio of air:"
Replace-Unique "= Specific heat ratio of air;" "= Specific heat ratio of air:"

# 13. "= Number of stages;" -> "= Number of stages:"
Replace-Unique "= Number of stages;" "= Number of stages:"

# 14. "= Conversion constant;" (C4, followed by "3.03") -> colon
Replace-Unique "= Conversion constant; 3.03" "= Conversion constant: 3.03"

# 15. "Compressor operating pressure;" -> "Compressor operating pressure:"
Replace-Unique "Compressor operating pressure;" "Compressor operating pressure:"

# 16. "= Compressor motor efficiency;" -> "= Compressor motor efficiency:"
Replace-Unique "= Compressor motor efficiency;" "= Compressor motor efficiency:"

# 17. "= Conversion constant;" (C5, followed by "0.746 kW/ HP") -> colon
Replace-Unique "= Conversion constant; 0.746" "= Conversion constant: 0.746"

# 18. "= Annual time during which leak occurs; " -> "...occurs: "
Replace-Unique "leak occurs; " "leak occurs: "

# 19. "...the facility peak demand; 1" -> "...peak demand: 1"
Replace-Unique "facility peak demand; 1" "facility peak demand: 1"

# 20. "= Conversion constant;" (C6, followed by "12 mos") -> colon
Replace-Unique "= Conversion constant; 12 mos" "= Conversion constant: 12 mos"
